# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary values ---
$ws.Range("E11").Value = 547925      # VALOR MORA
$ws.Range("C13").Value = 3           # Cant. Trabajadores
$ws.Range("F13").Value = 8           # Cant. Periodos

# --- Copy the "last row" special formatting (currently on row 25) up to row 23  ---
# before we overwrite / remove rows, so the new last data row (23) ends up with
# the same bottom-border style that the old last data row (25) had.
$ws.Range("B25:J25").Copy()
$ws.Range("B23:J23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Rewrite the worker/period detail table (rows 16-23) ---
# Row 16: GLEINER ENRIQUE CHARRIS COTE - periodo 1901
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "19562002"
$ws.Range("D16").Value = "GLEINER ENRIQUE CHARRIS COTE"
$ws.Range("E16").Value = "1901"
$ws.Range("F16").Value = 10800
$ws.Range("G16").Value = 900000

# Row 17: GLEINER ENRIQUE CHARRIS COTE - periodo 1902
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "19562002"
$ws.Range("D17").Value = "GLEINER ENRIQUE CHARRIS COTE"
$ws.Range("E17").Value = "1902"
$ws.Range("F17").Value = 36000
$ws.Range("G17").Value = 900000

# Row 18: GLEINER ENRIQUE CHARRIS COTE - periodo 1903
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "19562002"
$ws.Range("D18").Value = "GLEINER ENRIQUE CHARRIS COTE"
$ws.Range("E18").Value = "1903"
$ws.Range("F18").Value = 36000
$ws.Range("G18").Value = 900000

# Row 19: HECTOR ENRIQUE VILL SALAS - periodo 1907
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1047423245"
$ws.Range("D19").Value = "HECTOR ENRIQUE VILL SALAS"
$ws.Range("E19").Value = "1907"
$ws.Range("F19").Value = 33125
$ws.Range("G19").Value = 877803

# Row 20: GERSEY ENRIQUE MORALES FERRER - periodo 2008
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "72247754"
$ws.Range("D20").Value = "GERSEY ENRIQUE MORALES FERRER"
$ws.Range("E20").Value = "2008"
$ws.Range("F20").Value = 108000
$ws.Range("G20").Value = 2700000

# Row 21: GERSEY ENRIQUE MORALES FERRER - periodo 2009
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "72247754"
$ws.Range("D21").Value = "GERSEY ENRIQUE MORALES FERRER"
$ws.Range("E21").Value = "2009"
$ws.Range("F21").Value = 108000
$ws.Range("G21").Value = 2700000

# Row 22: GERSEY ENRIQUE MORALES FERRER - periodo 2010
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "72247754"
$ws.Range("D22").Value = "GERSEY ENRIQUE MORALES FERRER"
$ws.Range("E22").Value = "2010"
$ws.Range("F22").Value = 108000
$ws.Range("G22").Value = 2700000

# Row 23: GERSEY ENRIQUE MORALES FERRER - periodo 2011 (last data row, special style)
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "72247754"
$ws.Range("D23").Value = "GERSEY ENRIQUE MORALES FERRER"
$ws.Range("E23").Value = "2011"
$ws.Range("F23").Value = 108000
$ws.Range("G23").Value = 2700000

# Remove the now-unused last two rows of the old (longer) table.
$ws.Rows("24:25").Delete()
